$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two new columns (ownTeam, oppTeam) before the old
# "batsman" column (D), pushing batsman..sr from D..I to F..K.
$ws.Columns("D:E").Insert()

# Insert a new row above the existing data row for the new Sept 24 2020
# match, pushing the existing Sept 21 2020 row from row 2 down to row 3.
$ws.Rows(2).Insert()

# --- Header row (row 1): new column headers ---
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# --- New row 2: September 24 2020 match vs Kings XI Punjab ---
$ws.Range("A2").Value = " Dubai (DSC)"
$ws.Range("B2").Value = " September 24 2020"
$ws.Range("C2").Value = "Kings XI won by 97 runs"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Kings XI Punjab"
$ws.Range("F2").Value = "AB de Villiers "

# Numeric-looking values are stored as text on this sheet (same as the
# existing rows), so force the Text number format before writing them -
# otherwise Excel auto-converts them to real numbers.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "28"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "18"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "4"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "1"

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "155.55"

# --- Row 3 (previously row 2): fill in the new ownTeam / oppTeam columns ---
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Sunrisers Hyderabad"

Write-Host "done"
